# Update example input spreadsheets
#
# - transitions!C3: "time_dependent_gompertz" -> "time_dependent_weibull"
# - transitions!D13: new empty styled cell (font "Var(--jp-code-font-family)"
#   size 10, left/center aligned) -- pushes the sheet dimension out to row 13
# - specification!B1: max_iterations 500 -> 1000
# - specification!B2: time_horizon 10 -> 15
# - remove the now-unused "condensed_states" worksheet

$wb = $excel.ActiveWorkbook

$transitions = $wb.Worksheets.Item("transitions")
$transitions.Range("C3").Value = "time_dependent_weibull"

$d13 = $transitions.Range("D13")
$d13.Font.Name = "Var(--jp-code-font-family)"
$d13.Font.Size = 10
$d13.HorizontalAlignment = -4131
$d13.VerticalAlignment = -4108

$d3e3 = $transitions.Range("D3:E3")
$d3e3.WrapText = $false

$spec = $wb.Worksheets.Item("specification")
$spec.Range("B1").Value = 1000
$spec.Range("B2").Value = 15

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("condensed_states").Delete()
